$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scale the "value" column (D), rows 2-33, by 10000.
# Row 29 is blank and is skipped automatically since it has no numeric value.
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value2 = [Math]::Round($val * 10000, 6)
    }
}
